# Apply the LinuxForHealth re-branding + metadata refresh to the
# StructureDefinition-line-of-business workbook.
#
# Summary of the change (per the supplied OOXML diff):
#   Metadata sheet:
#     URL       : http://ibm.com/fhir/cdm/StructureDefinition/line-of-business
#                 -> http://linuxforhealth.org/fhir/cdm/StructureDefinition/line-of-business
#     Version   : 7.0.0 -> 8.0.0
#     Date      : 2022-09-08T16:11:15+00:00 -> 2022-11-10T16:00:46+00:00
#     Publisher : Alvearie Team -> LinuxForHealth Team
#   Elements sheet:
#     Fixed Value of Extension.url (row 5, "Q" column) mirrors the same URL
#     string, so it also needs to move to the linuxforhealth.org host.
#     The Constraint(s) cell for the root "Extension" row (row 2, column AI)
#     is cleared - that ele-1/ext-1 constraint note only belongs on the
#     Extension.extension row (row 4), which already carries it.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/line-of-business"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/line-of-business"
$elements.Range("AI2").Value = ""
